# Add "NA" values under the duplicate_image_filename column (column E)
# for the data rows 2 through 21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
